# Actualización automática 2025-06-30 15:05:09
# Applies newly recorded sales for asesor LOZANO MOLINA TITO and propagates
# the resulting totals/summary figures across the three report sheets.

$wb = $excel.ActiveWorkbook

$wsGrupo  = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": newly registered purchases for two clients
# of asesor LOZANO MOLINA TITO.
# ---------------------------------------------------------------------

# Row 14 - LINO TUMBACO VICENTE JAVIER
$wsGrupo.Range("E14").Value = 64.81999999999999   # FREGADEROS DE COCINA
$wsGrupo.Range("N14").Value = 762.16               # PUERTAS DE SEGURIDAD
$wsGrupo.Range("Q14").Value = 646.75               # PANELES PVC

# Row 18 - PAREDES ORTIZ MARIA INES
$wsGrupo.Range("E18").Value = 64.81999999999999   # FREGADEROS DE COCINA

# Row 29 - totals ("X de 27" counters) updated to reflect the two new
# non-zero entries above.
$wsGrupo.Range("E29").Value = "4 de 27"
$wsGrupo.Range("N29").Value = "1 de 27"
$wsGrupo.Range("Q29").Value = "1 de 27"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": junio (column F) sales updated for the same
# two clients, plus the column total.
# ---------------------------------------------------------------------

$wsMensual.Range("F14").Value = 1473.73
$wsMensual.Range("F18").Value = 64.81999999999999
$wsMensual.Range("F29").Value = 15470.5

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": VENTA / POR CUMPLIR / CUMPLIMIENTO
# figures recalculated for the affected groups and the TOTAL row.
# ---------------------------------------------------------------------

# Row 4 - FREGADEROS DE COCINA
$wsCumpl.Range("D4").Value = 2067.1
$wsCumpl.Range("E4").Value = -1816.468174579099
$wsCumpl.Range("F4").Value = 8.247555937992294

# Row 14 - PANELES PVC
$wsCumpl.Range("D14").Value = 1276.23
$wsCumpl.Range("E14").Value = -793.23
$wsCumpl.Range("F14").Value = 2.642298136645963

# Row 17 - PUERTAS DE SEGURIDAD
$wsCumpl.Range("D17").Value = 762.16
$wsCumpl.Range("E17").Value = -420.16
$wsCumpl.Range("F17").Value = 2.228538011695906

# Row 19 - TOTAL
$wsCumpl.Range("D19").Value = 22702.45
$wsCumpl.Range("E19").Value = 797.5509300503857
$wsCumpl.Range("F19").Value = 0.9660616638942116
